# Supplementary table: strain distribution across infection groups.
# "F1 hybrids" row counts grow (this cohort absorbed the rows that used to
# be broken out separately), and the "M. m. domesticus" / "M. m. musculus"
# breakout rows are removed entirely now that the design/immune-genes
# table folds them into the F1 hybrids figures.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 3 = "F1 hybrids" data row -> update its four numeric cells.
$row = $t.Rows.Item(3)
$row.Cells.Item(2).Range.Text = "38"
$row.Cells.Item(3).Range.Text = "34"
$row.Cells.Item(4).Range.Text = "38"
$row.Cells.Item(5).Range.Text = "110"

# Remove the "M. m. domesticus" and "M. m. musculus" rows (old rows 4 & 5).
# Deleting row 4 twice removes both, since row 5 shifts up to index 4 once
# the first delete happens.
$t.Rows.Item(4).Delete()
$t.Rows.Item(4).Delete()
